$wb = $excel.ActiveWorkbook

# --- Sheet: Neodymium ---
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030

$ws.Range("B2").Value = [double]"0"
$ws.Range("C2").Value = [double]"2.72967848032444E-06"
$ws.Range("D2").Value = [double]"0.006425129170340081"
$ws.Range("E2").Value = [double]"0.009981179422717998"

$ws.Range("B3").Value = [double]"3.676266793030401E-12"
$ws.Range("C3").Value = [double]"0.0001301849549774341"
$ws.Range("D3").Value = [double]"0.005650646918323966"
$ws.Range("E3").Value = [double]"0.008375274494933352"

$ws.Range("B4").Value = [double]"5.738750001259246E-14"
$ws.Range("C4").Value = [double]"0.0001176308004029253"
$ws.Range("D4").Value = [double]"0.00471063984337027"
$ws.Range("E4").Value = [double]"0.006754746205633251"

$ws.Range("C5").Value = [double]"2.655036418867232E-09"
$ws.Range("D5").Value = [double]"0.0002552879416942394"
$ws.Range("E5").Value = [double]"0.0005205069953698282"

# --- Sheet: Dysprosium ---
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030

# --- Sheet: Copper ---
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030

$ws.Range("B2").Value = [double]"6.274753108837452E-06"
$ws.Range("C2").Value = [double]"0.004777218731695482"
$ws.Range("D2").Value = [double]"0.5852873017284498"
$ws.Range("E2").Value = [double]"0.7342234181682109"

$ws.Range("B3").Value = [double]"4.265998591901206E-05"
$ws.Range("C3").Value = [double]"0.01726122059679833"
$ws.Range("D3").Value = [double]"0.4314739977014601"
$ws.Range("E3").Value = [double]"0.5292569039985252"

$ws.Range("B4").Value = [double]"0.0001264939349260171"
$ws.Range("C4").Value = [double]"0.004648601452327738"
$ws.Range("D4").Value = [double]"0.3635370532549198"
$ws.Range("E4").Value = [double]"0.4609895213897682"

$ws.Range("B5").Value = [double]"3.974642869529463E-05"
$ws.Range("C5").Value = [double]"0.01017152216720051"
$ws.Range("D5").Value = [double]"0.5292536818480308"
$ws.Range("E5").Value = [double]"0.6059072426603161"

# --- Sheet: Raw silicon ---
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030

$ws.Range("B2").Value = [double]"1.076099402252142E-06"
$ws.Range("C2").Value = [double]"7.453705115974282E-05"
$ws.Range("D2").Value = [double]"0.01585703061242156"
$ws.Range("E2").Value = [double]"0.01666371100707123"

$ws.Range("B3").Value = [double]"1.148298763755546E-06"
$ws.Range("C3").Value = [double]"0.0002477726001429727"
$ws.Range("D3").Value = [double]"0.008711650735289463"
$ws.Range("E3").Value = [double]"0.009403475543899196"

$ws.Range("B4").Value = [double]"7.353957242698485E-06"
$ws.Range("C4").Value = [double]"6.997433578934436E-05"
$ws.Range("D4").Value = [double]"0.009332746844457358"
$ws.Range("E4").Value = [double]"0.01175525920959753"

$ws.Range("B5").Value = [double]"3.950110178780338E-06"
$ws.Range("C5").Value = [double]"8.880409956795554E-05"
$ws.Range("D5").Value = [double]"0.01528869999002913"
$ws.Range("E5").Value = [double]"0.01410775496041511"

Write-Output "edits applied"
